$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in C1, matching the style of the existing header cells (A1/B1)
$ws.Range("C1").Value = "standard_error"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108

# Update column B values (uppskattat_antal_ripor) and add column C values (standard_error)
$ws.Range("B2").Value = 78.15710943778942
$ws.Range("C2").Value = 11.88448578118932

$ws.Range("B3").Value = 27.92424284047004
$ws.Range("C3").Value = 5.998082589560713

$ws.Range("B4").Value = 9.851052335388042
$ws.Range("C4").Value = 3.294424791765418

$ws.Range("B5").Value = 82.78195239821885
$ws.Range("C5").Value = 12.3916320879611

$ws.Range("B6").Value = 44.48862345013954
$ws.Range("C6").Value = 8.057825718361977

$ws.Range("B7").Value = 29.80150286335878
$ws.Range("C7").Value = 6.242865669599347

$ws.Range("B8").Value = 13.02618490629824
$ws.Range("C8").Value = 3.844283452493729

$ws.Range("B9").Value = 43.78245482394686
$ws.Range("C9").Value = 7.973622128637855

$ws.Range("B10").Value = 78.15710943778942
$ws.Range("C10").Value = 11.88448578118932

$ws.Range("B11").Value = 48.84819339861839
$ws.Range("C11").Value = 8.57265064267621
